$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing "Morocco" -> "Maroc" translation (D42 keeps same shared-string slot)
$ws.Range("D42").Value = "Maroc"

# New row 43: Kamel Mhalhel
$ws.Range("A43").Value = "Kamel"
$ws.Range("B43").Value = "Mhalhel"
$ws.Range("C43").Value = "Università degli studi di Messina"
$ws.Range("D43").Value = "Italie"
$ws.Range("E43").Value = "R35oeVAAAAAJ"
$ws.Range("F2").Copy()
$ws.Range("F43").PasteSpecial(-4122)
$ws.Range("F43").Value = "M"
$ws.Range("G43").Value = 1992
$ws.Range("H43").Value = "Médecine, Biologie et Sciences de la Santé"

# New row 44: Akram Zribi
$ws.Range("A44").Value = "Akram"
$ws.Range("B44").Value = "Zribi"
$ws.Range("C44").Value = "Université de Tunis El Manar"
$ws.Range("D44").Value = "Tunisie"
$ws.Range("E44").Value = "FAZ-BeAAAAAJ"
$ws.Range("F2").Copy()
$ws.Range("F44").PasteSpecial(-4122)
$ws.Range("F44").Value = "M"
$ws.Range("G44").Value = 1985
$ws.Range("H44").Value = "Chimie et Sciences des Matériaux"

# New row 45: Sondes Mechri
$ws.Range("A45").Value = "Sondes"
$ws.Range("B45").Value = "Mechri"
$ws.Range("C45").Value = "Université de Sfax"
$ws.Range("D45").Value = "Tunisie"
$ws.Range("E45").Value = "RN8eVNQAAAAJ"
$ws.Range("F3").Copy()
$ws.Range("F45").PasteSpecial(-4122)
$ws.Range("F45").Value = "F"
$ws.Range("G45").Value = 1991
$ws.Range("H45").Value = "Médecine, Biologie et Sciences de la Santé"

$ws.Range("H46").Select()
